$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text (e.g. "19.53"). Plain
# `.Value = "19.53"` would let Excel auto-convert it to a real number,
# which does not match the original inline-string cell content. Force
# those specific cells to Text format first so the assignment keeps the
# literal string, then restore the cell style to Normal/General so no
# stray number-format style lingers on the cell (values like
# "25.954.24" that contain two dots never parse as a number, so they
# don't need this treatment).
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"

# Updated "Price" (column D) values.
$ws.Range("D2").Value = "25.954.24"
$ws.Range("D3").Value = "1.638.57"
$ws.Range("D5").Value = "214.67"
$ws.Range("D10").Value = "19.53"
$ws.Range("D13").Value = "1.603.28"
$ws.Range("D17").Value = "25.981.67"
$ws.Range("D19").Value = "194.11"
$ws.Range("D24").Value = "143.76"
$ws.Range("D30").Value = "0.0494"
$ws.Range("D35").Value = "0.901"
$ws.Range("D36").Value = "1.129.88"
$ws.Range("D40").Value = "98.49"
$ws.Range("D41").Value = "5.42"
$ws.Range("D42").Value = "0.794"
$ws.Range("D47").Value = "7.73"
$ws.Range("D49").Value = "0.999"
$ws.Range("D50").Value = "0.0946"

# Restore default ("Normal") style on the cells we temporarily forced to
# text format, so the cell reverts to the General number format/style
# it started with (only the stored text changes, not the style index).
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(50, 4).Style = "Normal"

# Updated "Volume(1h)" (column E) percentage values. These already carry
# padding spaces and a trailing "%" so Excel never mistakes them for
# numbers, and can be set directly.
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("E23").Value = "  +3.73%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  -0.90%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("E51").Value = "  -0.62%  "
